$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename C1 from "Signal" to "Base Signal", add a new D1 column
#     "Joint 4 Signal" for the second slip-ring signal wire needed for final assembly. ---
$ws.Range("C1").Value = "Base Signal"
$ws.Range("D1").Value = "Joint 4 Signal"
$ws.Range("D1").Font.Bold = $true
# Touch E1 (kept blank) so the header row/dimension extends through column E, matching
# the extra spacer column next to the new "Joint 4 Signal" column.
$ws.Range("E1").Font.Bold = $true

# --- New "Joint 4 Signal" column (D) mirrors column C (the base signal wire) for most
#     channels; rows 8-11 differ because those channels pair a motor STP/DIR line on C
#     with a GND/24V return wire on D. ---
$ws.Range("C2").Value  = "GND";        $ws.Range("D2").Value  = "GND"
$ws.Range("C3").Value  = "GND";        $ws.Range("D3").Value  = "GND"
$ws.Range("C4").Value  = "GND";        $ws.Range("D4").Value  = "GND"
$ws.Range("C5").Value  = "GND";        $ws.Range("D5").Value  = "GND"
$ws.Range("C6").Value  = "GND";        $ws.Range("D6").Value  = "GND"
$ws.Range("C7").Value  = "COM";        $ws.Range("D7").Value  = "COM"
$ws.Range("C8").Value  = "Mot 2 STP";  $ws.Range("D8").Value  = "GND"
$ws.Range("C9").Value  = "Mot 2 DIR";  $ws.Range("D9").Value  = "GND"
$ws.Range("C10").Value = "Mot 3 STP";  $ws.Range("D10").Value = "24V"
$ws.Range("C11").Value = "Mot 3 DIR";  $ws.Range("D11").Value = "24V"
$ws.Range("C12").Value = "24V";        $ws.Range("D12").Value = "24V"
$ws.Range("C13").Value = "Mot 4 STP";  $ws.Range("D13").Value = "Mot 4 STP"
$ws.Range("C14").Value = "Mot 4 DIR";  $ws.Range("D14").Value = "Mot 4 DIR"
$ws.Range("C15").Value = "Mot 5 STP";  $ws.Range("D15").Value = "Mot 5 STP"
$ws.Range("C16").Value = "Mot 5 DIR";  $ws.Range("D16").Value = "Mot 5 DIR"
$ws.Range("C17").Value = "Mot 6 STP";  $ws.Range("D17").Value = "Mot 6 STP"
$ws.Range("C18").Value = "Mot 6 DIR";  $ws.Range("D18").Value = "Mot 6 DIR"
$ws.Range("C19").Value = "24V";        $ws.Range("D19").Value = "24V"
$ws.Range("C20").Value = "CAN L";      $ws.Range("D20").Value = "CAN L"
$ws.Range("C21").Value = "CAN H";      $ws.Range("D21").Value = "CAN H"
$ws.Range("C22").Value = "EN";         $ws.Range("D22").Value = "EN"
$ws.Range("C23").Value = "24V";        $ws.Range("D23").Value = "24V"
$ws.Range("C24").Value = "24V";        $ws.Range("D24").Value = "24V"
$ws.Range("C25").Value = "24V";        $ws.Range("D25").Value = "24V"

# --- Column C keeps the same "best fit" width behaviour as column A ---
$ws.Columns("C").ColumnWidth = $ws.Columns("A").ColumnWidth

# --- View tweaks made while finishing the sheet for rendering ---
$excel.ActiveWindow.Zoom = 190
$ws.Range("F8").Select()
